# Update "想去人数" (F) / "最低票价" (G) figures on the 展览 and 全部类型
# sheets to match the newly scraped numbers (gh-pages output generated at
# 456a3b4).

$wb = $excel.ActiveWorkbook

# Column F updates shared between the "展览" sheet (rows keyed by its own
# row numbers) and the "全部类型" sheet (offset by the extra 演出 row that
# sheet also carries). Column G only changes for one row (19.9 -> 60).

$sheetRowMap = @{
    "展览"   = @{
        7  = @{ F = 5211 }
        8  = @{ F = 170 }
        11 = @{ F = 67 }
        14 = @{ F = 12 }
        15 = @{ F = 6476 }
        19 = @{ F = 161 }
        21 = @{ F = 15457; G = 60 }
        22 = @{ F = 1535 }
        23 = @{ F = 289 }
        24 = @{ F = 146 }
        26 = @{ F = 11091 }
        27 = @{ F = 760 }
        28 = @{ F = 4340 }
        29 = @{ F = 248 }
        32 = @{ F = 306 }
    }
    "全部类型" = @{
        8  = @{ F = 5211 }
        9  = @{ F = 170 }
        13 = @{ F = 67 }
        17 = @{ F = 12 }
        18 = @{ F = 6476 }
        22 = @{ F = 161 }
        24 = @{ F = 15457; G = 60 }
        25 = @{ F = 1535 }
        26 = @{ F = 289 }
        27 = @{ F = 146 }
        29 = @{ F = 11091 }
        30 = @{ F = 760 }
        31 = @{ F = 4340 }
        32 = @{ F = 248 }
        35 = @{ F = 306 }
    }
}

foreach ($sheetName in $sheetRowMap.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetRowMap[$sheetName]
    foreach ($row in $rows.Keys) {
        $cols = $rows[$row]
        if ($cols.ContainsKey("F")) {
            $ws.Cells.Item($row, 6).Value = $cols["F"]
        }
        if ($cols.ContainsKey("G")) {
            $ws.Cells.Item($row, 7).Value = $cols["G"]
        }
    }
}
